$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 465, shifting the existing rows (465-479) down to (466-480).
$ws.Rows.Item(465).Insert()

# Populate the newly inserted row 465 with the new weekly data point.
$ws.Range("A465").Value = 10
$ws.Range("B465").Value = "Vega Modelo de Temuco"
$ws.Range("C465").Value = "La Araucanía"
$ws.Range("D465").Value = 45075
$ws.Range("D465").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E465").Value = 9
$ws.Range("F465").Value = 100112017
$ws.Range("G465").Value = "Apio"
$ws.Range("H465").Value = "Americana (o)"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 110
$ws.Range("K465").Value = 8000
$ws.Range("L465").Value = 8000
$ws.Range("M465").Value = 8000
$ws.Range("N465").Value = '$/docena de matas'
$ws.Range("O465").Value = "Provincia del Elquí"
$ws.Range("P465").Value = 1333
$ws.Range("Q465").Value = 6
$ws.Range("R465").Value = "Hortaliza"
